# Adds two new classification template sheets (Local Government District and
# Assembly Area) after the existing "Classifications_EQ - Template" sheet,
# mirroring the layout/format of that sheet, and wires up the workbook-level
# AutoFilter defined names for the new sheets.

$wb = $excel.ActiveWorkbook
$eq = $wb.Worksheets.Item("Classifications_EQ - Template")

# ---------------------------------------------------------------------------
# 1. Create "Classifications_LGD - Template" right after the EQ template.
# ---------------------------------------------------------------------------
$lgd = $wb.Worksheets.Add($null, $eq)
$lgd.Name = "Classifications_LGD - Template"

# ---------------------------------------------------------------------------
# 2. Create "Classifications_AA - Template" right after the LGD template.
# ---------------------------------------------------------------------------
$aa = $wb.Worksheets.Add($null, $lgd)
$aa.Name = "Classifications_AA - Template"

# ---------------------------------------------------------------------------
# Helper data - Local Government Districts (11)
# ---------------------------------------------------------------------------
$lgdCodes = @("N09000001","N09000002","N09000003","N09000004","N09000005","N09000006","N09000007","N09000008","N09000009","N09000010","N09000011")
$lgdNames = @("Antrim and Newtownabbey","Armagh, Banbridge and Craigavon","Belfast","Causeway Coast and Glens","Derry and Strabane","Fermanagh and Omagh","Lisburn and Castlereagh","Mid and East Antrim","Mid Ulster","Newry, Mourne and Down","North Down and Ards")

# Helper data - Assembly Areas (18)
$aaCodes = @("N06000001","N06000002","N06000003","N06000004","N06000005","N06000006","N06000007","N06000008","N06000009","N06000010","N06000011","N06000012","N06000013","N06000014","N06000015","N06000016","N06000017","N06000015")
$aaNames = @("Belfast East","Belfast North","Belfast South","Belfast West","East Antrim","East Londonderry","Fermanagh and South Tyrone","Foyle","Lagan Valley","Mid Ulster","Newry and Armagh","North Antrim","North Down","South Antrim","South Down","Strangford","Upper Bann","West Tyrone")

# ---------------------------------------------------------------------------
# 3. Populate "Classifications_LGD - Template"
# ---------------------------------------------------------------------------
$lgd.Range("A1").Value = "CODE"
$lgd.Range("B1").Value = "VALUE"

$lgd.Range("A2").Value = "N92000002"
$lgd.Range("B2").Value = "Northern Ireland"
$lgd.Range("C2").Value = "Group"

for ($i = 0; $i -lt $lgdNames.Count; $i++) {
    $r = 3 + $i
    $lgd.Cells.Item($r, 2).Value = $lgdNames[$i]
    $lgd.Cells.Item($r, 3).Value = "Local Government District"
}
for ($i = 0; $i -lt $lgdCodes.Count; $i++) {
    $r = 3 + $i
    $lgd.Cells.Item($r, 1).Value = $lgdCodes[$i]
}

# ---------------------------------------------------------------------------
# 4. Populate "Classifications_AA - Template"
# ---------------------------------------------------------------------------
$aa.Range("A1").Value = "CODE"
$aa.Range("B1").Value = "VALUE"

$aa.Range("A2").Value = "N92000002"
$aa.Range("B2").Value = "Northern Ireland"
$aa.Range("C2").Value = "Group"

for ($i = 0; $i -lt $aaNames.Count; $i++) {
    $r = 3 + $i
    $aa.Cells.Item($r, 2).Value = $aaNames[$i]
    $aa.Cells.Item($r, 3).Value = "Assembly Area"
}
for ($i = 0; $i -lt $aaCodes.Count; $i++) {
    $r = 3 + $i
    $aa.Cells.Item($r, 1).Value = $aaCodes[$i]
}

# ---------------------------------------------------------------------------
# 5. Match formatting of the EQ template: bold-free "VALUE" column header /
#    Northern Ireland row font, column widths, trailing blank styled row.
# ---------------------------------------------------------------------------
$eq.Range("B1").Copy()
$lgd.Range("B1").PasteSpecial(-4122) | Out-Null
$aa.Range("B1").PasteSpecial(-4122) | Out-Null

$eq.Range("B2").Copy()
$lgd.Range("B2").PasteSpecial(-4122) | Out-Null
$aa.Range("B2").PasteSpecial(-4122) | Out-Null

$eq.Range("B74").Copy()
$lgd.Range("B14").PasteSpecial(-4122) | Out-Null
$aa.Range("B21").PasteSpecial(-4122) | Out-Null
$aa.Range("B21").Value = $null

$eq.Columns.Item("A:C").AutoFit() | Out-Null
$lgd.Columns.Item("A:C").AutoFit() | Out-Null
$aa.Columns.Item("A:C").AutoFit() | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6. AutoFilter ranges on the new sheets + workbook-level _FilterDatabase
#    defined names (mirrors the existing EQ template setup).
# ---------------------------------------------------------------------------
$lgd.Range("A2:C13").AutoFilter(1) | Out-Null
$aa.Range("A2:C20").AutoFilter(1) | Out-Null

$lgd.Names.Add("_xlnm._FilterDatabase", "='Classifications_LGD - Template'!`$A`$2:`$C`$13") | Out-Null
$aa.Names.Add("_xlnm._FilterDatabase", "='Classifications_AA - Template'!`$A`$2:`$C`$20") | Out-Null

# ---------------------------------------------------------------------------
# 7. Leave the newly-added Assembly Area sheet as the active tab (matches the
#    activeTab bump recorded in the workbook view after this edit).
# ---------------------------------------------------------------------------
$lgd.Range("A1").Select() | Out-Null
$aa.Range("A1").Select() | Out-Null
$aa.Activate()

Write-Output "Added Classifications_LGD - Template and Classifications_AA - Template sheets"
